# Update the AMPICI row of Table 1 so it completely aligns with the PENICI
# row (penicillinase-labile penicillins, e.g. ampicillin, should be reported
# Resistant whenever PENICI is Resistant - CLSI VET01S 6th ed, Table 2C-1,
# comment 9; see also Table 1 of academic.oup.com/cid/article/58/9/1287).
#
# Also bumps the two Cochran-Armitage trend P-values (ERYTH and PENICI rows)
# that shifted slightly as a result of the corrected AMPICI counts.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- AMPICI row (row 5) ------------------------------------------------
# Columns 2-15: per-year "not susceptible % (n isolates tested)" values.
# Column 16: Cochran-Armitage trend P-value.

$ampiciRow = 5
$ampiciReplacements = @(
    @{ Col = 2;  Old = "54 (57)";    New = "75 (57)" },
    @{ Col = 3;  Old = "44 (57)";    New = "74 (57)" },
    @{ Col = 4;  Old = "56 (41)";    New = "88 (41)" },
    @{ Col = 5;  Old = "55 (67)";    New = "84 (67)" },
    @{ Col = 6;  Old = "59 (92)";    New = "84 (92)" },
    @{ Col = 7;  Old = "40 (99)";    New = "79 (99)" },
    @{ Col = 8;  Old = "44 (86)";    New = "80 (86)" },
    @{ Col = 9;  Old = "50 (135)";   New = "85 (135)" },
    @{ Col = 10; Old = "49 (107)";   New = "83 (107)" },
    @{ Col = 11; Old = "60 (90)";    New = "81 (90)" },
    @{ Col = 12; Old = "55 (184)";   New = "76 (184)" },
    @{ Col = 13; Old = "66 (214)";   New = "86 (214)" },
    @{ Col = 14; Old = "72 (181)";   New = "87 (181)" },
    @{ Col = 15; Old = "56 (1410)";  New = "82 (1410)" },
    @{ Col = 16; Old = "<0.001";     New = "0.168" }
)

foreach ($item in $ampiciReplacements) {
    $cell = $t.Cell($ampiciRow, $item.Col)
    if ($cell.Range.Text -like ($item.Old + "*")) {
        $cell.Range.Text = $item.New
    }
}

# --- ERYTH row (row 16): trend P-value 0.14 -> 0.15 --------------------
$erythCell = $t.Cell(16, 16)
if ($erythCell.Range.Text -like "0.14*") {
    $erythCell.Range.Text = "0.15"
}

# --- PENICI row (row 23): trend P-value 0.149 -> 0.159 -----------------
$peniciCell = $t.Cell(23, 16)
if ($peniciCell.Range.Text -like "0.149*") {
    $peniciCell.Range.Text = "0.159"
}
